$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously overflowed with placeholder rows (4-10). Trim it
# back down to the real sample data (rows 1-3) as part of adding the
# manual "download sheet" / "delete rows" buttons and guarding against an
# overflowing sheet.

# Row 2 gains a 4th ("D") sample column, stored as text like its
# existing siblings (A2:C2 are text, not numbers).
$ws.Range("D2").Value = "'4"
$ws.Range("D2").ClearFormats()

# Row 3 likewise gains its D column value.
$ws.Range("D3").Value = "'4"
$ws.Range("D3").ClearFormats()

# Remove the now-unneeded overflow rows 4 through 10.
$ws.Range("A4:D10").EntireRow.Delete()
